$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.617.52"
$ws.Range("E2").Value = "  +8.10%  "
$ws.Range("D3").Value = "3.390.66"
$ws.Range("E3").Value = "  +4.41%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'417.90"
$ws.Range("E5").Value = "  +5.98%  "
$ws.Range("D6").Value = "'116.11"
$ws.Range("E6").Value = "  +8.01%  "
$ws.Range("D7").Value = "'0.598"
$ws.Range("E7").Value = "  +5.97%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "'0.648"
$ws.Range("E9").Value = "  +5.34%  "
$ws.Range("D10").Value = "'40.90"
$ws.Range("E10").Value = "  +5.24%  "
$ws.Range("D11").Value = "'0.104"
$ws.Range("E11").Value = "  +6.26%  "
$ws.Range("E12").Value = "  +1.16%  "
$ws.Range("D13").Value = "3.916.81"
$ws.Range("E13").Value = "  +4.22%  "
$ws.Range("E14").Value = "  +4.80%  "
$ws.Range("D15").Value = "'20.07"
$ws.Range("E15").Value = "  +5.92%  "
$ws.Range("D16").Value = "3.363.84"
$ws.Range("E16").Value = "  +3.66%  "
$ws.Range("E17").Value = "  +2.25%  "
$ws.Range("D18").Value = "61.251.79"
$ws.Range("E18").Value = "  +7.84%  "
$ws.Range("D19").Value = "'10.84"
$ws.Range("E19").Value = "  +1.00%  "
$ws.Range("E20").Value = "  +2.89%  "
$ws.Range("D21").Value = "'0.0000114"
$ws.Range("E21").Value = "  +8.07%  "
$ws.Range("D22").Value = "'13.17"
$ws.Range("E22").Value = "  +1.36%  "
$ws.Range("D23").Value = "'305.92"
$ws.Range("E23").Value = "  +2.86%  "
$ws.Range("D24").Value = "'76.06"
$ws.Range("E24").Value = "  +3.41%  "
$ws.Range("E25").Value = "  +2.94%  "
$ws.Range("D26").Value = "'29.08"
$ws.Range("E26").Value = "  +4.15%  "
$ws.Range("D27").Value = "'4.49"
$ws.Range("E27").Value = "  +2.32%  "
$ws.Range("D28").Value = "'8.00"
$ws.Range("E28").Value = "  +3.58%  "
$ws.Range("D29").Value = "'7.71"
$ws.Range("E29").Value = "  +7.03%  "
$ws.Range("E30").Value = "  +6.36%  "
$ws.Range("E31").Value = "  +6.59%  "
$ws.Range("E32").Value = "  +23.80%  "
$ws.Range("D33").Value = "'11.57"
$ws.Range("E33").Value = "  +5.64%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").Value = "'39.96"
$ws.Range("E35").Value = "  +7.66%  "
$ws.Range("D36").Value = "'0.0510"
$ws.Range("E36").Value = "  +5.54%  "
$ws.Range("D37").Value = "'52.67"
$ws.Range("E37").Value = "  +2.06%  "
$ws.Range("D38").Value = "'3.14"
$ws.Range("E38").Value = "  +2.50%  "
$ws.Range("D39").Value = "'0.997"
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("E40").Value = "  -2.04%  "
$ws.Range("D41").Value = "'137.38"
$ws.Range("E41").Value = "  +2.44%  "
$ws.Range("E42").Value = "  +3.90%  "
$ws.Range("D43").Value = "'0.295"
$ws.Range("E43").Value = "  +4.36%  "
$ws.Range("E44").Value = "  +2.22%  "
$ws.Range("D45").Value = "'4.02"
$ws.Range("E45").Value = "  +2.16%  "
$ws.Range("D46").Value = "'17.10"
$ws.Range("E46").Value = "  +0.64%  "
$ws.Range("D47").Value = "'22.77"
$ws.Range("E47").Value = "  +4.11%  "
$ws.Range("D48").Value = "'2.27"
$ws.Range("E48").Value = "  +10.01%  "
$ws.Range("D49").Value = "2.176.87"
$ws.Range("E49").Value = "  +1.57%  "
$ws.Range("E50").Value = "  +1.77%  "
$ws.Range("D51").Value = "'2.00"
$ws.Range("E51").Value = "  -1.50%  "
